$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, shifting the existing rows 198-201 down to 199-202.
$ws.Rows(198).Insert()

# Populate the newly inserted row 198 with the new weekly data point.
$ws.Range("A198").Value = 7
$ws.Range("B198").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C198").Value = 'Ñuble'
$ws.Range("D198").Value = 44595
$ws.Range("E198").Value = 16
$ws.Range("F198").Value = 100112043
$ws.Range("G198").Value = 'Pepino ensalada'
$ws.Range("H198").Value = 'Sin especificar'
$ws.Range("I198").Value = 'Primera'
$ws.Range("J198").Value = 80
$ws.Range("K198").Value = 9500
$ws.Range("L198").Value = 10000
$ws.Range("M198").Value = 9750
$ws.Range("N198").Value = '$/caja 80 unidades'
$ws.Range("O198").Value = 'Región del Maule'
$ws.Range("P198").Value = 122
$ws.Range("Q198").Value = 80
$ws.Range("R198").Value = 'Hortaliza'
